# Daily attendance processing - 2026-01-05 23:02:16
# Normalize the "Recorded By" column (G) so entries that contain
# "System" alongside the recorder's email are listed with the email
# first, e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$colG = $ws.Range("G1:G" + $usedRange.Rows.Count)

$colG.Replace("System, dnasr281@gmail.com", "dnasr281@gmail.com, System")
